$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91 (shifts old rows 91-189 down to 92-190)
$ws.Rows(91).Insert()

# Populate the newly inserted row 91 with the new record
$ws.Range("A91").Value = 4
$ws.Range("B91").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C91").Value = "Los Lagos"
$ws.Range("D91").Value = 44589
$ws.Range("E91").Value = 10
$ws.Range("F91").Value = "Fruta"
$ws.Range("G91").Value = 100103
$ws.Range("H91").Value = "Frutos de hueso (carozo)"
$ws.Range("I91").Value = 100103004
$ws.Range("J91").Value = "Durazno"
$ws.Range("K91").Value = "Carson"
$ws.Range("L91").Value = "Primera"
$ws.Range("M91").Value = 400
$ws.Range("N91").Value = 17000
$ws.Range("O91").Value = 18000
$ws.Range("P91").Value = 17500
$ws.Range("Q91").Value = "$/caja 15 kilos empedrada"
$ws.Range("R91").Value = "Región de O'Higgins"
$ws.Range("S91").Value = 1167
$ws.Range("T91").Value = 15
